$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.53611
$ws.Range("F2").Value = 0.20778
$ws.Range("I2").Value = 0.06054999999999999
$ws.Range("K2").Value = 1.99466
$ws.Range("M2").Value = 1.88371
$ws.Range("N2").Value = 0.00124
